# Applies the "feat: add 2022-Q1 data" change:
#  - the existing "总计" sheet (rId2 / sheetId 2) is renamed to "2022-Q1" and
#    repopulated with the new fund-holding detail rows
#  - a fresh "总计" sheet is appended right after it, re-stating the summary
#    table with a new row for 2022-Q1 (the old 2021-Q4 row shifts down)

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1. Turn the current "总计" sheet into the new "2022-Q1" detail sheet.
# ---------------------------------------------------------------------
$q1 = $wb.Worksheets.Item("总计")
$q1.Name = "2022-Q1"

# Clear out the old 总计 content (B1:D2) before writing the new layout.
$q1.Cells.Clear()

# Header row (bold, centered, boxed - same look as the existing header row
# on the "2021-Q4" sheet).
$headerRange = $q1.Range("B1:H1")
$headerRange.Value = ""
$headerRange.Font.Bold = $true
$headerRange.HorizontalAlignment = -4108
$headerRange.VerticalAlignment = -4160
$headerRange.Borders.LineStyle = 1

$q1.Range("B1").Value = "基金代码"
$q1.Range("C1").Value = "基金名称"
$q1.Range("D1").Value = "基金规模"
$q1.Range("E1").Value = "股票总仓位"
$q1.Range("F1").Value = "仓位占比"
$q1.Range("G1").Value = "持有市值(亿元)"
$q1.Range("H1").Value = "仓位排名"

# Data rows -- columns B-G are stored as text (matches the source data
# export), column A (index) and H (rank) are numeric.
$rows = @(
    @{ A = 0; B = "011052"; C = "鹏华弘裕一年持有期混合A"; D = "2.92"; E = "24.56"; F = "1.98"; G = "0.0578"; H = 4 },
    @{ A = 1; B = "001190"; C = "鹏华弘润灵活配置混合 - A"; D = "3.65"; E = "23.00"; F = "1.21"; G = "0.0442"; H = 4 },
    @{ A = 2; B = "011053"; C = "鹏华弘裕一年持有期混合C"; D = "0.20"; E = "24.56"; F = "1.98"; G = "0.0040"; H = 4 },
    @{ A = 3; B = "001191"; C = "鹏华弘润灵活配置混合 - C"; D = "0.25"; E = "23.00"; F = "1.21"; G = "0.0030"; H = 4 }
)

# Pre-format B2:G5 as text so the numeric-looking strings ("2.92", "24.56",
# ...) are stored verbatim instead of being coerced into numbers.
$q1.Range("B2:G5").NumberFormat = "@"

$indexRange = $q1.Range("A2:A5")
$indexRange.Font.Bold = $true
$indexRange.HorizontalAlignment = -4108
$indexRange.VerticalAlignment = -4160
$indexRange.Borders.LineStyle = 1

$r = 2
foreach ($row in $rows) {
    $q1.Cells.Item($r, 1).Value = $row.A
    $q1.Cells.Item($r, 2).Value = $row.B
    $q1.Cells.Item($r, 3).Value = $row.C
    $q1.Cells.Item($r, 4).Value = $row.D
    $q1.Cells.Item($r, 5).Value = $row.E
    $q1.Cells.Item($r, 6).Value = $row.F
    $q1.Cells.Item($r, 7).Value = $row.G
    $q1.Cells.Item($r, 8).Value = $row.H

    $r = $r + 1
}

# ---------------------------------------------------------------------
# 2. Re-create the "总计" summary sheet right after "2022-Q1".
# ---------------------------------------------------------------------
$total = $wb.Worksheets.Add($null, $q1)
$total.Name = "总计"

$totalHeader = $total.Range("B1:D1")
$totalHeader.Value = ""
$totalHeader.Font.Bold = $true
$totalHeader.HorizontalAlignment = -4108
$totalHeader.VerticalAlignment = -4160
$totalHeader.Borders.LineStyle = 1

$total.Range("B1").Value = "日期"
$total.Range("C1").Value = "持有数量(只)"
$total.Range("D1").Value = "持有市值(亿元)"

$totalIndex = $total.Range("A2:A3")
$totalIndex.Font.Bold = $true
$totalIndex.HorizontalAlignment = -4108
$totalIndex.VerticalAlignment = -4160
$totalIndex.Borders.LineStyle = 1

$total.Cells.Item(2, 1).Value = 0
$total.Cells.Item(2, 2).Value = "2022-Q1"
$total.Cells.Item(2, 3).Value = 4
$total.Cells.Item(2, 4).Value = 0.11

$total.Cells.Item(3, 1).Value = 1
$total.Cells.Item(3, 2).Value = "2021-Q4"
$total.Cells.Item(3, 3).Value = 1
$total.Cells.Item(3, 4).Value = 0

# Keep the originally active sheet/selection in place (no tab switch was
# part of this change).
$wb.Worksheets.Item("2021-Q4").Activate() | Out-Null
$wb.Worksheets.Item("2021-Q4").Range("A1").Select() | Out-Null

